$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4177.597448727024
$ws.Range("C3").Value = 4177.597448727024
$ws.Range("C4").Value = 4022.994289493579
$ws.Range("C5").Value = 4021.268557446408
$ws.Range("C6").Value = 4021.268557446408
$ws.Range("C7").Value = 4008.519560781017
$ws.Range("C8").Value = 3961.743297273461
$ws.Range("C9").Value = 3961.743297273461
$ws.Range("C10").Value = 3961.743297273461
$ws.Range("C11").Value = 3961.743297273461
$ws.Range("C12").Value = 3961.743297273461
